# Applies the data row permutation described by the diff:
#   new row 2 = old row 3
#   new row 3 = old row 5
#   new row 4 = old row 4 (unchanged)
#   new row 5 = old row 2
# Only columns D and J..Q vary between these rows; we rewrite that full
# range for rows 2, 3 and 5 using the values captured from the original
# (pre-edit) rows.
#
# Note: this runtime's .Value getter is unreliable, so .Value2 is used
# consistently for both reading and writing cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (rows 2, 3, 5) before overwriting anything,
# since rows 2 and 3 get overwritten and row 5 depends on the original row 2.
$origRow2 = @{
    D = $ws.Range("D2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    N = $ws.Range("N2").Value2
    O = $ws.Range("O2").Value2
    P = $ws.Range("P2").Value2
    Q = $ws.Range("Q2").Value2
}

$origRow3 = @{
    D = $ws.Range("D3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
    Q = $ws.Range("Q3").Value2
}

$origRow5 = @{
    D = $ws.Range("D5").Value2
    J = $ws.Range("J5").Value2
    K = $ws.Range("K5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    Q = $ws.Range("Q5").Value2
}

function Set-DataRow($rowNum, $data) {
    $ws.Range("D$rowNum").Value2 = $data.D
    $ws.Range("J$rowNum").Value2 = $data.J
    $ws.Range("K$rowNum").Value2 = $data.K
    $ws.Range("L$rowNum").Value2 = $data.L
    $ws.Range("M$rowNum").Value2 = $data.M
    $ws.Range("N$rowNum").Value2 = $data.N
    $ws.Range("O$rowNum").Value2 = $data.O
    $ws.Range("P$rowNum").Value2 = $data.P
    $ws.Range("Q$rowNum").Value2 = $data.Q
}

# New row 2 = old row 3
Set-DataRow 2 $origRow3

# New row 3 = old row 5
Set-DataRow 3 $origRow5

# New row 5 = old row 2
Set-DataRow 5 $origRow2

$wb.Save()
